# Applies the 31-10-2023 scraper update to the Serbia Super Liga 2023-2024 sheet:
#  1. Rows 2 and 3 had their match data swapped back to the correct fixtures.
#  2. Rows 77 and 78 had their match data swapped back to the correct fixtures.
#  3. Eight new fixtures (rows 88-95) were appended at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [object[]]$Values
    )
    $cols = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22)
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($Row, $cols[$i]).Value = $Values[$i]
    }
}

# --- Fix rows 2 and 3 (swap the two fixtures) ---
Set-Row 2 @(1, "serbia", "super-liga", "2023-2024", 45136.78819444445, "Sp. Subotica", 2, "IMT Novi Beograd", 1, `
    1.83, "28/07/2023 16:12", 1.93, "29/07/2023 18:45", 3.22, "28/07/2023 16:12", 3.28, "29/07/2023 18:45", `
    3.85, "28/07/2023 16:12", 4.06, "29/07/2023 18:45", `
    "https://www.betexplorer.com/football/serbia/super-liga/spartak-subotica-imt-novi-beograd/nwxR1PLi/")

Set-Row 3 @(2, "serbia", "super-liga", "2023-2024", 45136.78819444445, "Cukaricki", 2, "Radnicki Nis", 0, `
    1.71, "28/07/2023 07:12", 1.55, "29/07/2023 18:47", 3.34, "28/07/2023 07:12", 3.85, "29/07/2023 18:47", `
    4.3, "28/07/2023 07:12", 6, "29/07/2023 18:47", `
    "https://www.betexplorer.com/football/serbia/super-liga/cukaricki-radnicki-nis/xSUwaN5A/")

# --- Fix rows 77 and 78 (swap the two fixtures) ---
Set-Row 77 @(76, "serbia", "super-liga", "2023-2024", 45207.72916666666, "TSC", 0, "Sp. Subotica", 2, `
    1.4, "07/10/2023 04:43", 1.4, "08/10/2023 14:55", 4.25, "07/10/2023 04:43", 4.71, "08/10/2023 17:29", `
    6.3, "07/10/2023 04:43", 7.08, "08/10/2023 17:29", `
    "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-spartak-subotica/ADMQ0wxi/")

Set-Row 78 @(77, "serbia", "super-liga", "2023-2024", 45207.72916666666, "Cukaricki", 1, "Vozdovac", 1, `
    1.46, "07/10/2023 04:43", 1.46, "08/10/2023 17:29", 4.05, "07/10/2023 04:43", 4.12, "08/10/2023 17:29", `
    5.75, "07/10/2023 04:43", 7.03, "08/10/2023 17:29", `
    "https://www.betexplorer.com/football/serbia/super-liga/cukaricki-fk-vozdovac/WMNM1JMp/")

# --- Append the 8 new fixtures (rows 88-95), copying formatting from the last existing row first ---
$ws.Range("A87:V87").Copy()
$ws.Range("A88:V95").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Set-Row 88 @(87, "serbia", "super-liga", "2023-2024", 45226.75, "Javor", 0, "Napredak", 2, `
    1.85, "26/10/2023 06:12", 2.08, "27/10/2023 17:58", 3.14, "26/10/2023 06:12", 2.91, "27/10/2023 17:58", `
    3.91, "26/10/2023 06:12", 4.06, "27/10/2023 17:58", `
    "https://www.betexplorer.com/football/serbia/super-liga/javor-napredak/l0BkjSZ6/")

Set-Row 89 @(88, "serbia", "super-liga", "2023-2024", 45227.625, "Radnik", 0, "Zeleznicar Pancevo", 1, `
    1.96, "27/10/2023 03:12", 2.45, "28/10/2023 14:52", 3.15, "27/10/2023 03:12", 2.67, "28/10/2023 14:52", `
    3.47, "27/10/2023 03:12", 3.47, "28/10/2023 14:52", `
    "https://www.betexplorer.com/football/serbia/super-liga/radnik-surdulica-zeleznicar-pancevo/UelxEVRD/")

Set-Row 90 @(89, "serbia", "super-liga", "2023-2024", 45227.77083333334, "IMT Novi Beograd", 1, "Crvena zvezda", 2, `
    8.15, "27/10/2023 06:42", 24.2, "28/10/2023 18:29", 5.6, "27/10/2023 06:42", 9.53, "28/10/2023 18:29", `
    1.23, "27/10/2023 06:42", 1.09, "28/10/2023 18:21", `
    "https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-crvena-zvezda/SjAgknkD/")

Set-Row 91 @(90, "serbia", "super-liga", "2023-2024", 45227.77083333334, "Radnicki Nis", 1, "Sp. Subotica", 1, `
    1.7, "27/10/2023 06:42", 1.66, "28/10/2023 18:23", 3.42, "27/10/2023 06:42", 3.65, "28/10/2023 18:23", `
    4.23, "27/10/2023 06:42", 5.08, "28/10/2023 18:23", `
    "https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-spartak-subotica/2qDshl5f/")

Set-Row 92 @(91, "serbia", "super-liga", "2023-2024", 45228.54166666666, "Radnicki 1923", 4, "Mladost", 1, `
    1.71, "28/10/2023 02:13", 1.83, "29/10/2023 12:50", 3.41, "28/10/2023 02:13", 3.53, "29/10/2023 12:50", `
    4.37, "28/10/2023 02:13", 4.14, "29/10/2023 12:20", `
    "https://www.betexplorer.com/football/serbia/super-liga/radnicki-1923-mladost-lucani/tzCoi8K0/")

Set-Row 93 @(92, "serbia", "super-liga", "2023-2024", 45228.6875, "Cukaricki", 2, "Vojvodina", 2, `
    1.97, "28/10/2023 02:42", 2.12, "29/10/2023 15:56", 3.25, "28/10/2023 02:42", 3.18, "29/10/2023 15:56", `
    3.33, "28/10/2023 02:42", 3.51, "29/10/2023 16:28", `
    "https://www.betexplorer.com/football/serbia/super-liga/cukaricki-vojvodina/lxX6bjSQ/")

Set-Row 94 @(93, "serbia", "super-liga", "2023-2024", 45228.77083333334, "Partizan", 3, "Vozdovac", 0, `
    1.21, "28/10/2023 02:42", 1.2, "29/10/2023 18:26", 5.57, "28/10/2023 02:42", 6.34, "29/10/2023 18:28", `
    9.05, "28/10/2023 02:42", 11.9, "29/10/2023 18:28", `
    "https://www.betexplorer.com/football/serbia/super-liga/partizan-fk-vozdovac/KKNxgUkl/")

Set-Row 95 @(94, "serbia", "super-liga", "2023-2024", 45228.79166666666, "TSC", 1, "Novi Pazar", 1, `
    1.41, "28/10/2023 02:42", 1.36, "29/10/2023 18:52", 4.1, "28/10/2023 02:42", 4.69, "29/10/2023 18:52", `
    6.09, "28/10/2023 02:42", 8.38, "29/10/2023 18:52", `
    "https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-novi-pazar/6HJYgAzr/")

Write-Output "Applied update: swapped rows 2/3, 77/78, appended rows 88-95"
